# Update the code for Practitioner
# Get API for the Practitioner is completed
#
# - RGID "RG-34565541998716" -> "RG-20002"
# - Provider value on row 2 becomes the Practitioner id "P-388412033222"
# - New row 3 added with the Facility id "F-468464031024" (same RGID)
# - Header row gets a border added (to match the bordered data rows)
# - Columns auto-sized to the new, shorter content
# - Selection / view state updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------
$ws.Range("A2").Value = "RG-20002"
$ws.Range("B2").Value = "P-388412033222"

$ws.Range("A3").Value = "RG-20002"
$ws.Range("B3").Value = "F-468464031024"

# give the new row the same bordered style already used by row 2
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header formatting ------------------------------------------------
# add the thin border around the header cells (keeps existing bold/yellow fill)
$ws.Range("A1:B1").Borders.LineStyle = 1

# --- Column widths (autofit to the new, shorter values) ---------------
$ws.Columns("A:B").AutoFit()

# --- View / selection state --------------------------------------------
[void]$ws.Range("F10").Select()

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.Orientation = 1
